# Modulo de salida completado
# Applies the "salida" (withdrawal) edit to the warehouse request/delivery form.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header: update the "DIA" (day) value ---
$ws.Range("I3").Value = "17"

# --- Row 9: first existing line item is changed to a different article ---
$ws.Range("A9").Value = "21601"
$ws.Range("B9").Value = "Cloralex 1lt"
$ws.Range("G9").Value = "Pieza"
$ws.Range("H9").Value = 7.333333333333333

# --- Row 10: a brand new second line item is added ---
$ws.Range("A10").Value = "5000"
$ws.Range("B10").Value = "Silla de escritorio"
$ws.Range("F10").Value = 1
$ws.Range("G10").Value = "Pieza"
$ws.Range("H10").Value = 100
$ws.Range("J10").Value = "Ninguna"

# --- Signature block: new names for the responsible people ---
$ws.Range("B26").Value = "Jesus Alberto Calderón García"
$ws.Range("F26").Value = "Tila del Carmen Mendoza Olan"
$ws.Range("J26").Value = "Antonio Espinosa Correa"
